$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 159, pushing the existing rows 159:165 down to 160:166
$ws.Rows.Item(159).Insert()

# Fill in the new weekly record at row 159
$ws.Cells.Item(159, 1).Value = 7
$ws.Cells.Item(159, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(159, 3).Value = "Ñuble"
$ws.Cells.Item(159, 4).Value = 44568
$ws.Cells.Item(159, 5).Value = 16
$ws.Cells.Item(159, 6).Value = 100112017
$ws.Cells.Item(159, 7).Value = "Apio"
$ws.Cells.Item(159, 8).Value = "Americana (o)"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 80
$ws.Cells.Item(159, 11).Value = 8000
$ws.Cells.Item(159, 12).Value = 8500
$ws.Cells.Item(159, 13).Value = 8250
$ws.Cells.Item(159, 14).Value = "`$/docena de matas"
$ws.Cells.Item(159, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(159, 16).Value = 1375
$ws.Cells.Item(159, 17).Value = 6
$ws.Cells.Item(159, 18).Value = "Hortaliza"
